$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - MAE
$ws.Range("B2").Value = 0.85
$ws.Range("C2").Value = 0.804
$ws.Range("D2").Value = 0.735
$ws.Range("E2").Value = 1.28
$ws.Range("F2").Value = 1.08

# Row 3 - MSE
$ws.Range("B3").Value = 1.484
$ws.Range("C3").Value = 1.178
$ws.Range("D3").Value = 1.18
$ws.Range("E3").Value = 6.739
$ws.Range("F3").Value = 2.975

# Row 4 - mean Y-Test
$ws.Range("B4").Value = 18.214
$ws.Range("C4").Value = 15.308
$ws.Range("D4").Value = 12.948
$ws.Range("E4").Value = 30.588
$ws.Range("F4").Value = 18.064

# Row 5 - mean Y-predicted
$ws.Range("B5").Value = 18.105
$ws.Range("C5").Value = 15.281
$ws.Range("D5").Value = 13.107
$ws.Range("E5").Value = 30.21
$ws.Range("F5").Value = 17.994

# Row 6 - R2
$ws.Range("B6").Value = 0.874
$ws.Range("C6").Value = 0.931
$ws.Range("D6").Value = 0.782
$ws.Range("E6").Value = 0.844
$ws.Range("F6").Value = 0.866
